$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46003, 11203.9679631759, 9680.254896972, 16560.26, 6676.89426313305, -8.46295166228932),
    @(46004, 4170.35788785527, 6586.51898992352, 11232.26, 6551.23394232305, 79.3955388436073),
    @(46005, 4107.91390949668, 6757.78246436953, 11232.26, 6664.19524347906, 91.2382378270246),
    @(46006, 11281.1952980293, 10555.7159208018, 11232.26, 7455.21027635104, 282.444424881369),
    @(46007, 11009.641624268, 10090.3041094634, 11232.26, 7266.50807737139, 255.189674451449),
    @(46008, 9746.03784374469, 9088.0715137262, 11232.26, 8042.27723337764, 245.753697795993),
    @(46009, 9746.03784374469, 9201.45129547426, 11232.26, 8042.27723337764, 250.477855368829),
    @(46010, 9746.03784374469, 8623.57393502478, 11232.26, 8042.27723337764, 226.399632016768),
    @(46011, 8560.57524882407, 8319.44127193012, 11232.26, 7650.34018821925, 197.396727506224),
    @(46012, 8469.04674334209, 8431.12133603587, 11232.26, 7641.94217211311, 201.700146172874),
    @(46013, 9566.18303866457, 9312.70103230373, 11232.26, 8113.15183942236, 258.066369655254),
    @(46014, 9566.18303866457, 9812.85930553035, 11232.26, 8113.15183942236, 278.906297706363),
    @(46015, 9566.18303866457, 9906.79025109728, 11232.26, 8113.15183942236, 282.820087104985),
    @(46016, 8471.80799626025, 9364.98722035868, 11232.26, 7768.77536493624, 245.895941053955)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
